# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap Huelva (row 53) and Huesca (row 54) in the "Ciudad" (A) and
# "Casos activos" (C) columns - the other columns (B, D, E) are identical
# between the two rows so no change is needed there.
$ws.Range("A53").Value = "Huesca"
$ws.Range("C53").Value = 0

$ws.Range("A54").Value = "Huelva"
$ws.Range("C54").Value = 72

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 22:16"
